# Add 2022-Q4 data (commit: "feat: add 2022-Q4 data")
#
# - Inserts a new "2022-Q4" worksheet right after the "总计" summary sheet
#   (ahead of the existing "2022-Q3" sheet), populated with the quarter's
#   fund-holding table.
# - Inserts a matching row at the top of the "总计" sheet's data so the
#   new quarter shows up in the rolled-up history, renumbering the index
#   column for the rows that shift down.

$wb = $excel.ActiveWorkbook
$summarySheet = $wb.Worksheets.Item(1)

# A style-"2" cell (bold, centered, thin border) to clone formatting from -
# every sheet in this workbook stamps its header row / index column with
# this exact look, so we reuse it via copy/paste-format instead of
# re-deriving the look property-by-property (which would mint new style
# indexes).
$styleTemplate = $summarySheet.Cells.Item(1, 2)

function Copy-HeaderStyle($targetCell) {
    $styleTemplate.Copy() | Out-Null
    $targetCell.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q4Sheet.Name = "2022-Q4"

# Header row (row 1): A1 stays empty, B1..H1 carry the column titles.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q4Sheet.Cells.Item(1, $i + 2)
    Copy-HeaderStyle $cell
    $cell.Value = $headers[$i]
}

# Data rows 2..11 - one row per fund holding.
# Columns: code, name, fund size, stock position, position ratio, holding
# value (亿元), position rank.
$rows = @(
    @("011164", "富国兴远优选12个月持有期混合A", "37.17", "85.15", "3.31", "1.2303", 9),
    @("011165", "富国兴远优选12个月持有期混合C", "15.84", "85.15", "3.31", "0.5243", 9),
    @("001186", "富国文体健康股票A",             "10.99", "83.45", "4.41", "0.4847", 7),
    @("001150", "融通互联网传媒灵活配置混合",       "8.57",  "90.42", "3.37", "0.2888", 4),
    @("011830", "富国天恒混合A",                 "3.31",  "87.21", "4.14", "0.1370", 9),
    @("004809", "新疆前海联合润丰灵活配置混合A",     "1.23",  "90.88", "4.06", "0.0499", 5),
    @("011125", "富国文体健康股票C",             "0.87",  "83.45", "4.41", "0.0384", 7),
    @("004890", "中邮健康文娱灵活配置混合",         "0.42",  "92.60", "6.63", "0.0278", 3),
    @("011831", "富国天恒混合C",                 "0.04",  "87.21", "4.14", "0.0017", 9),
    @("005935", "新疆前海联合润丰灵活配置混合C",     "0.03",  "90.88", "4.06", "0.0012", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]

    # Column A: numeric 0-based row index, styled like the header.
    $idxCell = $q4Sheet.Cells.Item($rowNum, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $r

    # Columns B..G are text in the source data (fund code keeps leading
    # zeros, numeric-looking figures are kept as text) - force text
    # format before assigning so COM doesn't coerce them to numbers, then
    # drop back to the Normal style so no stray numFmt/style lingers
    # (matches the source, which leaves these cells unstyled).
    for ($c = 2; $c -le 7; $c++) {
        $cell = $q4Sheet.Cells.Item($rowNum, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $data[$c - 2]
        $cell.Style = "Normal"
    }

    # Column H (仓位排名) is a plain, unstyled number.
    $q4Sheet.Cells.Item($rowNum, 8).Value = $data[6]
}

# ---------------------------------------------------------------------
# 2) Update "总计": insert a new row 2 for 2022-Q4, push the existing
#    quarters down, and renumber column A (0,1,2,...).
# ---------------------------------------------------------------------
$summarySheet.Cells.Item(2, 1).EntireRow.Insert()

$summaryRows = @(
    @("2022-Q4", 10, 2.78),
    @("2022-Q3", 8, 0.19),
    @("2021-Q4", 7, 1.54),
    @("2021-Q3", 5, 1.57),
    @("2021-Q2", 3, 0.15),
    @("2021-Q1", 2, 0.05)
)

for ($r = 0; $r -lt $summaryRows.Length; $r++) {
    $rowNum = $r + 2
    $data = $summaryRows[$r]

    $idxCell = $summarySheet.Cells.Item($rowNum, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $r

    # EntireRow.Insert() on row 2 copied the header row's bold styling
    # down into the new row - strip it back to Normal (unstyled, like
    # every other data row in this column) before/after writing values.
    $bCell = $summarySheet.Cells.Item($rowNum, 2)
    $bCell.Value = $data[0]
    $bCell.Style = "Normal"

    $cCell = $summarySheet.Cells.Item($rowNum, 3)
    $cCell.Value = $data[1]
    $cCell.Style = "Normal"

    $dCell = $summarySheet.Cells.Item($rowNum, 4)
    $dCell.Value = $data[2]
    $dCell.Style = "Normal"
}
